$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 2 (GLD) ---
$ws.Range("D2").Value = 386.36
$ws.Range("E2").Value = 68.59999999999999
$ws.Range("F2").Value = -0.39
$ws.Range("G2").Value = 50
$ws.Range("H2").Value = 73
$ws.Range("I2").Value = 80
$ws.Range("K2").Value = 62.2
$ws.Range("N2").Value = 50.68470204858703

# --- Row 3 (GC=F) ---
$ws.Range("D3").Value = 4249.1
$ws.Range("E3").Value = 72.3
$ws.Range("F3").Value = 0.73
$ws.Range("I3").Value = 76
$ws.Range("J3").Value = 80
$ws.Range("K3").Value = 60.6
$ws.Range("M3").Value = "📈 매수 관찰 구간입니다."
$ws.Range("N3").Value = 50.68470204858703

# --- Row 4 (NEM) ---
$ws.Range("D4").Value = 90.29000000000001
$ws.Range("E4").Value = 56.3
$ws.Range("F4").Value = -0.48
$ws.Range("N4").Value = 50.68470204858703
